# redmine # 9229 Calibration sheets added/changed for GP05MOAS gliders
# GL276, GL361-GL365, GL453, GL523, GL525, GL537, PG514, PG515.
#
# This particular workbook (GP05MOAS-GL361) had its Moorings mooring-cruise
# info updated (cruise vessel + recovery note) together with corrected
# anchor-launch / recover dates, and the Asset_Cal_Info sheet had one
# calibration coefficient value corrected.

$wb = $excel.ActiveWorkbook

# ---- Moorings sheet -------------------------------------------------
$moorings = $wb.Worksheets.Item("Moorings")

# Anchor Launch Date: 2014-02-09 -> 2014-02-20
$moorings.Range("D2").Value2 = 41690
# Recover Date: 2014-08-17 -> 2014-05-04
$moorings.Range("F2").Value2 = 41763
# Cruise Number: "Melville 130" -> "CCGS Tully"
$moorings.Range("J2").Value2 = "CCGS Tully"
# Notes: (blank) -> "Lost at sea"
$moorings.Range("K2").Value2 = "Lost at sea"

# ---- Asset_Cal_Info sheet --------------------------------------------
$assetCal = $wb.Worksheets.Item("Asset_Cal_Info")

# CC_angular_resolution calibration coefficient value: 1.13 -> 1.096
$assetCal.Range("F6").Value2 = 1.096
